# Updates the cryptocurrency price/volume snapshot data in Sheet1.
# Applies the latest scraped Price (column D) and Volume(1h) (column E)
# values for each coin row, matching the GitHub Actions data refresh.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '64.310.86'
$ws.Range("E2").Value = '  +0.04%  '
$ws.Range("D3").Value = '3.494.49'
$ws.Range("E3").Value = '  +0.00%  '
$ws.Range("E4").Value = '  +0.08%  '
$ws.Range("D5").Value = "'589.05"
$ws.Range("E5").Value = '  +0.37%  '
$ws.Range("D6").Value = "'134.15"
$ws.Range("E6").Value = '  -0.07%  '
$ws.Range("D8").Value = "'0.487"
$ws.Range("E8").Value = '  +0.24%  '
$ws.Range("D9").Value = "'7.61"
$ws.Range("E9").Value = '  +6.33%  '
$ws.Range("D10").Value = "'0.125"
$ws.Range("E10").Value = '  +0.13%  '
$ws.Range("D11").Value = "'0.389"
$ws.Range("E11").Value = '  +3.29%  '
$ws.Range("D12").Value = '4.089.72'
$ws.Range("E12").Value = '  -0.02%  '
$ws.Range("E13").Value = '  +0.78%  '
$ws.Range("E14").Value = '  -0.05%  '
$ws.Range("D15").Value = '3.504.37'
$ws.Range("E15").Value = '  +0.25%  '
$ws.Range("D16").Value = '64.266.58'
$ws.Range("E16").Value = '  -0.12%  '
$ws.Range("E17").Value = '  +0.56%  '
$ws.Range("D18").Value = "'10.02"
$ws.Range("E18").Value = '  +1.97%  '
$ws.Range("D19").Value = "'5.78"
$ws.Range("E19").Value = '  +0.79%  '
$ws.Range("D20").Value = "'13.56"
$ws.Range("E20").Value = '  -0.21%  '
$ws.Range("D21").Value = "'388.08"
$ws.Range("E21").Value = '  +0.09%  '
$ws.Range("D22").Value = "'0.581"
$ws.Range("E22").Value = '  +3.00%  '
$ws.Range("D23").Value = '3.632.62'
$ws.Range("E23").Value = '  -0.05%  '
$ws.Range("D24").Value = "'74.18"
$ws.Range("E24").Value = '  -0.31%  '
$ws.Range("E25").Value = '  +0.19%  '
$ws.Range("E26").Value = '  -0.57%  '
$ws.Range("E27").Value = '  +1.92%  '
$ws.Range("E28").Value = '  +0.04%  '
$ws.Range("D29").Value = "'7.37"
$ws.Range("E29").Value = '  -0.30%  '
$ws.Range("E30").Value = '  -1.66%  '
$ws.Range("E31").Value = '  +1.07%  '
$ws.Range("D32").Value = "'8.17"
$ws.Range("E32").Value = '  -1.21%  '
$ws.Range("E33").Value = '  +3.79%  '
$ws.Range("D34").Value = '3.523.51'
$ws.Range("E34").Value = '  +0.20%  '
$ws.Range("E35").Value = '  +0.01%  '
$ws.Range("D36").Value = "'23.33"
$ws.Range("E36").Value = '  -0.58%  '
$ws.Range("D37").Value = "'5.33"
$ws.Range("E37").Value = '  +1.89%  '
$ws.Range("D38").Value = "'6.94"
$ws.Range("E38").Value = '  +0.82%  '
$ws.Range("E39").Value = '  +0.49%  '
$ws.Range("D40").Value = "'165.48"
$ws.Range("E40").Value = '  +2.07%  '
$ws.Range("D41").Value = "'0.0787"
$ws.Range("E41").Value = '  +0.57%  '
$ws.Range("E42").Value = '  +0.30%  '
$ws.Range("E43").Value = '  -0.02%  '
$ws.Range("D44").Value = "'4.43"
$ws.Range("E44").Value = '  +0.66%  '
$ws.Range("D45").Value = "'24.52"
$ws.Range("E45").Value = '  -3.67%  '
$ws.Range("E46").Value = '  -0.35%  '
$ws.Range("E47").Value = '  -0.74%  '
$ws.Range("D48").Value = "'6.83"
$ws.Range("E48").Value = '  +1.29%  '
$ws.Range("D49").Value = '2.406.09'
$ws.Range("E49").Value = '  -2.82%  '
$ws.Range("D50").Value = "'0.920"
$ws.Range("E50").Value = '  +2.32%  '
$ws.Range("D51").Value = "'0.0259"
$ws.Range("E51").Value = '  -0.50%  '
